# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
#
# On the "Metadata" worksheet:
#   - B4 (the value for the "Name" property) is set to "PaysnationaliteVs"
#   - B8 (the value for the "Date" property) is updated to the new timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "PaysnationaliteVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
